$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E7").Value = 16.17130000000001
$ws.Range("A9").Value = -21.9496
$ws.Range("E12").Value = 18.00660000000003
$ws.Range("E14").Value = 16.78260000000001
$ws.Range("A18").Value = -22.20330000000001
$ws.Range("A20").Value = -21.09159999999996
$ws.Range("E26").Value = 16.2392
$ws.Range("A27").Value = -22.1372
$ws.Range("E27").Value = 16.53539999999999
$ws.Range("E29").Value = 17.09410000000001
$ws.Range("A35").Value = -20.89599999999997
$ws.Range("E37").Value = 16.70310000000001
$ws.Range("E38").Value = 16.6619
$ws.Range("E51").Value = 17.20550000000001
$ws.Range("E52").Value = 16.92750000000001
$ws.Range("E55").Value = 16.56190000000001
$ws.Range("A69").Value = -21.65149999999997
$ws.Range("E69").Value = 17.36160000000003
$ws.Range("E70").Value = 18.02830000000002
$ws.Range("A76").Value = -19.81779999999999
$ws.Range("A78").Value = -19.89409999999999
$ws.Range("E81").Value = 16.74399999999999
$ws.Range("A82").Value = -21.83560000000001
$ws.Range("A83").Value = -21.87159999999999
$ws.Range("E83").Value = 16.4928
$ws.Range("A93").Value = -21.11989999999999
$ws.Range("E102").Value = 16.81949999999999
